# Auto-generated Excel COM-interop script
# Applies cell-value updates described by the commit diff across all 8 sheets.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 1804.2142
$ws.Range("I9").Value = 1917.7778
$ws.Range("J9").Value = 1599.8
$ws.Range("K9").Value = 1917.7778
$ws.Range("L9").Value = 1599.8
$ws.Range("M9").Value = -1748.7778
$ws.Range("N9").Value = -1937.8
# Row 15
$ws.Range("H15").Value = 2032.5156
$ws.Range("I15").Value = 2032.5156
$ws.Range("K15").Value = 6097.5468
$ws.Range("M15").Value = -5928.5468
# Row 88
$ws.Range("H88").Value = 1127.9375
$ws.Range("I88").Value = 894.5
$ws.Range("J88").Value = 1205.75
$ws.Range("K88").Value = 894.5
$ws.Range("L88").Value = 1205.75
$ws.Range("M88").Value = -488.5
$ws.Range("N88").Value = -2017.75
# Row 91
$ws.Range("H91").Value = 1127.9375
$ws.Range("I91").Value = 894.5
$ws.Range("J91").Value = 1205.75
$ws.Range("K91").Value = 894.5
$ws.Range("L91").Value = 1205.75
$ws.Range("M91").Value = 509.5
$ws.Range("N91").Value = -4013.75
# Row 105
$ws.Range("H105").Value = 59999
$ws.Range("J105").Value = 59999
$ws.Range("L105").Value = 59999
$ws.Range("N105").Value = -66987
# Row 132
$ws.Range("H132").Value = 8337.932000000001
$ws.Range("I132").Value = 6453.1797
$ws.Range("J132").Value = 23039
$ws.Range("K132").Value = 19359.5391
$ws.Range("L132").Value = 69117
$ws.Range("M132").Value = -16829.5391
$ws.Range("N132").Value = -74177
# Row 135
$ws.Range("H135").Value = 2120
$ws.Range("I135").Value = 2390.6667
$ws.Range("K135").Value = 21516.0003
$ws.Range("M135").Value = -18981.0003
# Row 137
$ws.Range("H137").Value = 9313.666999999999
$ws.Range("I137").Value = 1684.4445
$ws.Range("K137").Value = 5053.333500000001
$ws.Range("M137").Value = -2503.333500000001
# Row 141
$ws.Range("H141").Value = 2170.4
$ws.Range("I141").Value = 1856
$ws.Range("K141").Value = 5568
$ws.Range("M141").Value = -388

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 26
$ws.Range("H26").Value = 499
$ws.Range("I26").Value = 499
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 499
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -169
$ws.Range("N26").ClearContents()
# Row 122
$ws.Range("H122").Value = 6074
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 36
$ws.Range("H36").Value = 6908
$ws.Range("I36").Value = 1124.75
$ws.Range("K36").Value = 1124.75
$ws.Range("M36").Value = -590.75
# Row 86
$ws.Range("H86").Value = 1929.2609
$ws.Range("I86").Value = 1963.4
$ws.Range("K86").Value = 1963.4
$ws.Range("M86").Value = -840.4000000000001
# Row 89
$ws.Range("H89").Value = 1929.2609
$ws.Range("I89").Value = 1963.4
$ws.Range("K89").Value = 9817
$ws.Range("M89").Value = -4201

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 29
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()
# Row 31
$ws.Range("H31").Value = 18293.436
$ws.Range("I31").Value = 6522.091
$ws.Range("J31").Value = 33526.94
$ws.Range("K31").Value = 6522.091
$ws.Range("L31").Value = 33526.94
$ws.Range("M31").Value = -6227.091
$ws.Range("N31").Value = -34116.94
# Row 34
$ws.Range("H34").Value = 18293.436
$ws.Range("I34").Value = 6522.091
$ws.Range("J34").Value = 33526.94
$ws.Range("K34").Value = 6522.091
$ws.Range("L34").Value = 33526.94
$ws.Range("M34").Value = -6320.091
$ws.Range("N34").Value = -33930.94
# Row 112
$ws.Range("H112").Value = 39999.332
$ws.Range("J112").Value = 39999.332
$ws.Range("L112").Value = 39999.332
$ws.Range("N112").Value = -42953.332

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 32
$ws.Range("H32").Value = 14288244
$ws.Range("J32").Value = 16667952
$ws.Range("L32").Value = 50003856
$ws.Range("N32").Value = -50004422
# Row 34
$ws.Range("H34").Value = 357862.9
$ws.Range("I34").Value = 1527.125
$ws.Range("J34").Value = 500397.2
$ws.Range("K34").Value = 4581.375
$ws.Range("L34").Value = 1501191.6
$ws.Range("M34").Value = -4497.375
$ws.Range("N34").Value = -1501359.6
# Row 39
$ws.Range("H39").Value = 3044.2727
$ws.Range("J39").Value = 2798.2856
$ws.Range("L39").Value = 8394.856800000001
$ws.Range("N39").Value = -8982.856800000001
# Row 55
$ws.Range("H55").Value = 5004
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 5004
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 15012
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -15366
# Row 122
$ws.Range("H122").Value = 12660130
$ws.Range("I122").Value = 31145848
$ws.Range("J122").Value = 2577011
$ws.Range("K122").Value = 280312632
$ws.Range("L122").Value = 23193099
$ws.Range("M122").Value = -280310182
$ws.Range("N122").Value = -23197999
# Row 125
$ws.Range("H125").Value = 6999
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
# Row 132
$ws.Range("H132").Value = 2100
$ws.Range("I132").Value = 1608.3334
$ws.Range("K132").Value = 14475.0006
$ws.Range("M132").Value = -11945.0006
# Row 134
$ws.Range("H134").Value = 6361.7026
$ws.Range("I134").Value = 3063.8333
$ws.Range("J134").Value = 7000
$ws.Range("K134").Value = 9191.499899999999
$ws.Range("L134").Value = 21000
$ws.Range("M134").Value = -4121.499899999999
$ws.Range("N134").Value = -31140

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 31
$ws.Range("H31").Value = 1152.091
$ws.Range("I31").Value = 1152.091
$ws.Range("K31").Value = 1152.091
$ws.Range("M31").Value = -860.0909999999999
# Row 37
$ws.Range("H37").Value = 1152.091
$ws.Range("I37").Value = 1152.091
$ws.Range("K37").Value = 1152.091
$ws.Range("M37").Value = -875.0909999999999
# Row 43
$ws.Range("H43").Value = 642504
$ws.Range("I43").Value = 642504
$ws.Range("K43").Value = 642504
$ws.Range("M43").Value = -642353
# Row 46
$ws.Range("H46").Value = 5431
$ws.Range("I46").Value = 960.25
$ws.Range("K46").Value = 960.25
$ws.Range("M46").Value = -804.25
# Row 80
$ws.Range("H80").Value = 6551.727
$ws.Range("J80").Value = 12473.454
$ws.Range("L80").Value = 12473.454
$ws.Range("N80").Value = -14469.454
# Row 83
$ws.Range("H83").Value = 6551.727
$ws.Range("J83").Value = 12473.454
$ws.Range("L83").Value = 62367.27
$ws.Range("N83").Value = -72351.26999999999
# Row 102
$ws.Range("H102").Value = 1846.7826
$ws.Range("I102").Value = 1628.15
$ws.Range("J102").Value = 3304.3333
$ws.Range("K102").Value = 1628.15
$ws.Range("L102").Value = 3304.3333
$ws.Range("M102").Value = -6.150000000000091
$ws.Range("N102").Value = -6548.3333
# Row 132
$ws.Range("H132").Value = 11302.3
$ws.Range("I132").Value = 8268.704
$ws.Range("K132").Value = 24806.112
$ws.Range("M132").Value = -22276.112

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 3972821.8
$ws.Range("I68").Value = 2349.25
$ws.Range("J68").Value = 5107242.5
$ws.Range("K68").Value = 2349.25
$ws.Range("L68").Value = 5107242.5
$ws.Range("M68").Value = -1600.25
$ws.Range("N68").Value = -5108740.5
# Row 71
$ws.Range("H71").Value = 3972821.8
$ws.Range("I71").Value = 2349.25
$ws.Range("J71").Value = 5107242.5
$ws.Range("K71").Value = 11746.25
$ws.Range("L71").Value = 25536212.5
$ws.Range("M71").Value = -8002.25
$ws.Range("N71").Value = -25543700.5
# Row 74
$ws.Range("H74").Value = 52183.332
$ws.Range("I74").Value = 45787.5
$ws.Range("K74").Value = 45787.5
$ws.Range("M74").Value = -44789.5
# Row 77
$ws.Range("H77").Value = 52183.332
$ws.Range("I77").Value = 45787.5
$ws.Range("K77").Value = 137362.5
$ws.Range("M77").Value = -132370.5
# Row 82
$ws.Range("H82").Value = 6961.4614
$ws.Range("J82").Value = 8999.75
$ws.Range("L82").Value = 8999.75
$ws.Range("N82").Value = -9721.75
# Row 85
$ws.Range("H85").Value = 6961.4614
$ws.Range("J85").Value = 8999.75
$ws.Range("L85").Value = 8999.75
$ws.Range("N85").Value = -11495.75

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 110
$ws.Range("H110").Value = 62820
$ws.Range("I110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("M110").ClearContents()
# Row 132
$ws.Range("H132").Value = 29051.25
$ws.Range("I132").Value = 4100
$ws.Range("J132").Value = 54002.5
$ws.Range("K132").Value = 12300
$ws.Range("L132").Value = 162007.5
$ws.Range("M132").Value = -9770
$ws.Range("N132").Value = -167067.5

